# Actualización automática 2025-08-15 16:30:08
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Narrow column F (CUMPLIMIENTO) from OOXML width 26 to 23
$ws.Range("F1").ColumnWidth = 22.17

# Update VENTA (D) for PORCELANATO row and recompute dependent columns
$ws.Range("D3").Value = 2848.66
$ws.Range("E3").Value = 12622.8993
$ws.Range("F3").Value = 0.184122359276353

# Update TOTAL row (D4/E4/F4) to reflect the new VENTA total
$ws.Range("D4").Value = 73834.88
$ws.Range("E4").Value = -58363.3207
$ws.Range("F4").Value = 4.772297256424568
